$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.076.03'
$ws.Range("E2").Value = '  -2.09%  '
$ws.Range("D3").Value = '2.260.02'
$ws.Range("E3").Value = '  -3.53%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '298.48'
$ws.Range("E5").Value = '  -2.74%  '
$ws.Range("D6").Value = '94.14'
$ws.Range("E6").Value = '  -6.97%  '
$ws.Range("D7").Value = '0.497'
$ws.Range("E7").Value = '  -2.67%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").Value = '0.493'
$ws.Range("E9").Value = '  -3.78%  '
$ws.Range("D10").Value = '33.01'
$ws.Range("E10").Value = '  -5.69%  '
$ws.Range("D11").Value = '0.0788'
$ws.Range("E11").Value = '  -1.73%  '
$ws.Range("D12").Value = '47.57'
$ws.Range("E12").Value = '  -8.61%  '
$ws.Range("D13").Value = '0.113'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").Value = '6.67'
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").Value = '2.608.17'
$ws.Range("E15").Value = '  -3.64%  '
$ws.Range("D16").Value = '15.26'
$ws.Range("E16").Value = '  -3.99%  '
$ws.Range("D17").Value = '2.261.08'
$ws.Range("E17").Value = '  -4.69%  '
$ws.Range("D18").Value = '0.775'
$ws.Range("E18").Value = '  -4.67%  '
$ws.Range("D19").Value = '42.045.24'
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").Value = '0.0₃0893'
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("D21").Value = '6.01'
$ws.Range("E21").Value = '  -3.59%  '
$ws.Range("D22").Value = '11.38'
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("D23").Value = '66.60'
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("D24").Value = '233.29'
$ws.Range("E24").Value = '  -1.59%  '
$ws.Range("D25").Value = '1.92'
$ws.Range("E25").Value = '  -5.20%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = '2.45'
$ws.Range("E27").Value = '  -4.31%  '
$ws.Range("D28").Value = '23.73'
$ws.Range("E28").Value = '  -7.35%  '
$ws.Range("D29").Value = '2.16'
$ws.Range("E29").Value = '  -7.12%  '
$ws.Range("D30").Value = '167.08'
$ws.Range("E30").Value = '  +4.37%  '
$ws.Range("D31").Value = '33.59'
$ws.Range("E31").Value = '  -4.51%  '
$ws.Range("D32").Value = '9.04'
$ws.Range("E32").Value = '  -3.52%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("D34").Value = '4.94'
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("D35").Value = '2.34'
$ws.Range("E35").Value = '  -5.55%  '
$ws.Range("D36").Value = '0.0694'
$ws.Range("E36").Value = '  -4.77%  '
$ws.Range("D37").Value = '4.38'
$ws.Range("E37").Value = '  -6.93%  '
$ws.Range("D38").Value = '2.79'
$ws.Range("E38").Value = '  -6.16%  '
$ws.Range("D39").Value = '15.97'
$ws.Range("E39").Value = '  -8.52%  '
$ws.Range("D40").Value = '0.0990'
$ws.Range("E40").Value = '  -3.88%  '
$ws.Range("D41").Value = '0.109'
$ws.Range("E41").Value = '  -3.44%  '
$ws.Range("D42").Value = '1.71'
$ws.Range("E42").Value = '  -8.65%  '
$ws.Range("D43").Value = '2.40'
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("D44").Value = '1.946.99'
$ws.Range("E44").Value = '  -3.65%  '
$ws.Range("D45").Value = '0.0278'
$ws.Range("E45").Value = '  -2.82%  '
$ws.Range("D46").Value = '17.39'
$ws.Range("E46").Value = '  -7.17%  '
$ws.Range("D47").Value = '9.55'
$ws.Range("E47").Value = '  -7.27%  '
$ws.Range("D48").Value = '2.79'
$ws.Range("E48").Value = '  -5.70%  '
$ws.Range("D49").Value = '2.81'
$ws.Range("E49").Value = '  -3.14%  '
$ws.Range("D50").Value = '2.483.10'
$ws.Range("E50").Value = '  -3.22%  '
$ws.Range("D51").Value = '52.19'
$ws.Range("E51").Value = '  -7.53%  '
